$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the time-range labels in column C
$ws.Range("C2").Value = "9:05-9:10"
$ws.Range("C3").Value = "9:10-9:15"

# Move the active cell selection from C14 to C11
$ws.Range("C11").Select()
